# Finish the "submitting PO based invoice manually" test data section:
#  - POBasedInvoice: refresh invoice numbers / base amounts / IGST, and make the
#    Quantity column mirror the (text) Base Amount value instead of a flat 1.
#  - Switch the active tab over to BADashboardPage.
#  - Widen BADashboardPage's state column now that it holds "TRIPURA" data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("POBasedInvoice")

# New invoice number / base amount / IGST per row. Quantity (J) is set to the
# same text as Base Amount (B), mirroring how the sheet was filled in by hand.
$rows = @(
    @{ Row = 2;  Invoice = "TESTINV24257"; Amount = "7"; Igst = "1.26" },
    @{ Row = 3;  Invoice = "TESTINV70287"; Amount = "8"; Igst = "1.44" },
    @{ Row = 4;  Invoice = "TESTINV07707"; Amount = "9"; Igst = "1.62" },
    @{ Row = 5;  Invoice = "TESTINV49593"; Amount = "9"; Igst = "1.62" },
    @{ Row = 6;  Invoice = "TESTINV87128"; Amount = "3"; Igst = "0.54" },
    @{ Row = 7;  Invoice = "TESTINV79234"; Amount = "8"; Igst = "1.44" },
    @{ Row = 8;  Invoice = "TESTINV09626"; Amount = "3"; Igst = "0.54" },
    @{ Row = 9;  Invoice = "TESTINV34765"; Amount = "7"; Igst = "1.26" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Invoice number (column A) is plain text already, safe to assign directly.
    $ws.Range("A$rowNum").Value = $r.Invoice

    # Base Amount (B), IGST (C) and Quantity (J) are stored as TEXT even though
    # they look numeric, so stage them as text via a helper formula cell and
    # paste-special (values) across - a direct .Value assignment would coerce
    # the numeric-looking string into a real number.
    $ws.Range("Z1").Formula = "=""" + $r.Amount + """"
    $ws.Range("Z1").Copy()
    $ws.Range("B$rowNum").PasteSpecial(-4163)

    $ws.Range("Z1").Formula = "=""" + $r.Igst + """"
    $ws.Range("Z1").Copy()
    $ws.Range("C$rowNum").PasteSpecial(-4163)

    $ws.Range("Z1").Formula = "=""" + $r.Amount + """"
    $ws.Range("Z1").Copy()
    $ws.Range("J$rowNum").PasteSpecial(-4163)
}

# Clean up the scratch cell used for staging text values.
$ws.Range("Z1").Clear()

# Move the active tab from POBasedInvoice to BADashboardPage.
$dash = $wb.Worksheets.Item("BADashboardPage")

# The "To state" sample value moves on from MANIPUR to TRIPURA.
$dash.Range("B2").Value = "TRIPURA"

$dash.Activate()

# Widen column B (To state) now that it shows the longer "TRIPURA" value.
$dash.Columns.Item(2).ColumnWidth = 13.83
